$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44449, 4, 10, 56.91843588138198),
    @(44450, 1, 9, 51.22659229324378),
    @(44451, 2, 9, 51.22659229324378),
    @(44452, 2, 9, 51.22659229324378),
    @(44453, 0, 9, 51.22659229324378),
    @(44454, 0, 9, 51.22659229324378),
    @(44455, 7, 16, 91.06949741021117),
    @(44456, 2, 14, 79.68581023393477),
    @(44457, 3, 16, 91.06949741021117),
    @(44458, 4, 18, 102.4531845864876),
    @(44459, 2, 18, 102.4531845864876)
)

$startRow = 375

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy formatting (style incl. number format) from the row above so new
    # rows look exactly like the existing data block.
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
